# Generate Report for Archive
#
# 1. Every "Status" cell currently showing "Ready for handoff" moves on to
#    "In Translation" (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4).
# 2. The two narrower "Status"-ish columns (Overview columns E & F, and the
#    "Status" column C on the per-language report sheets) get re-sized down
#    from their old width to the new, slightly narrower width.

$wb = $excel.ActiveWorkbook

# ---- 1. Update the "Ready for handoff" -> "In Translation" text ----------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F4").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C4").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C4").Value = "In Translation"

# ---- 2. Narrow the relevant columns ---------------------------------------
# ColumnWidth is expressed in characters; 12.5 is the closest request that
# lands the underlying pixel grid on the target column width.
$overview.Range("E1:F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
